$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet gains a new blank column before the existing
# "Late" column (N), pushing "Late" -> O and "Outstanding" -> Q, with a new
# blank column P in between (mirrors the old N/O blank-spacer pattern).
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# The "Repayment Schedule" tab becomes the active / selected sheet (it was
# "Transactions" before), with a new selected cell on that sheet.
$wsSchedule.Activate()
$wsSchedule.Range("L17").Select()
